$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume strings stay as Text (matching source formatting),
# mirroring the original workbook where columns D and E are text, not numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.304.85'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.20%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.606.37'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.25%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '546.16'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.49%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.08'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.38%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.20%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.88%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.37%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.334'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.31%  '

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.37%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.066.71'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.32%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '59.250.69'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.19%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.55'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.41%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.659.83'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.82%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000133'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.34%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '343.36'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.96%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.35'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.14%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.12'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.21%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.40'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.56%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.01%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.44'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.31%  '

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.96%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.407'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.14%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.75%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.20'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.63%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.00%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0737'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.54%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +8.61%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.81'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.77%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.75'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.39%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.10'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.07%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.98'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.13%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '37.09'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.94%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.08%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.24%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.834'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.68%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.813'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.40%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.51%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '277.37'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.38%  '

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.16%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.597'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.91%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.75'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.01%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0956'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.16%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.20%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.945.27'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.36%  '

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.69%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.32'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.80%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.51'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.52%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.87'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.95%  '
